$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty cells in row 2
$ws.Range("C2").Value = "ЕУР"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.01
$ws.Range("G2").Value = "test"
$ws.Range("H2").Value = "https://www.stb.com.mk/media/4548/%D1%82%D0%B0%D0%B1%D0%B5%D0%BB%D0%B8-%D0%B7%D0%B0-%D0%BA%D1%80%D0%B5%D0%B4%D0%B8%D1%82%D0%B8-%D0%B8-%D0%B4%D0%B5%D0%BF%D0%BE%D0%B7%D0%B8%D1%82%D0%B8-%D1%81%D0%BE-%D0%B4%D0%BE%D0%BF%D0%BE%D0%BB%D0%BD%D1%83%D0%B2%D0%B0nj%D0%B0-01-05-2025.xlsx"

# Move the active selection to C2 (matching the author's final cursor position)
[void]$ws.Range("C2").Select()
